# Add a new row (row 6) for phone 79174442 with 0 points, matching the
# existing rows' shape: phone stored as text, birthday left blank, points
# stored as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A6: phone number, entered as text (like the other "no birthday" rows would
# be if re-typed) rather than a plain number. Force text interpretation via
# NumberFormat so the numeric-looking string isn't coerced to a Number, then
# drop back to the Normal style so no stray formatting is left on the cell.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Formula = "79174442"
$ws.Range("A6").Style = "Normal"

# B6: birthday unknown/blank, same as B4/B5.
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Style = "Normal"

# C6: total_points starts at 0.
$ws.Range("C6").Value = 0
